# Update New Orleans xlsx file:
#   1. Insert a new "State" column (value "Louisiana") into the hotel_info
#      sheet, right after "Hotel_Name" and before "City".
#   2. Reorder the worksheet tabs so "review_info" comes before "hotel_info".
#
# NOTE: the column insert is performed BEFORE the sheet reorder so the
# $hotelWs handle still lines up with the worksheet we mean to edit.

$wb = $excel.ActiveWorkbook

$hotelWs = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# Insert a new column C ("State") in hotel_info, shifting City (and
# everything after it) one column to the right.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# Move review_info so it sits before hotel_info in the tab order.
$reviewWs.Move($hotelWs)
